$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "DatabaseName" column before the existing "Username" column (D),
# pushing Username -> E and Password -> F. Inherits formatting from the
# pushed column, matching the per-row look already on the sheet.
$ws.Columns("D").Insert()

# Header row
$ws.Range("D1").Value = "DatabaseName"

# Helper: write a digit-only string as literal TEXT (not auto-converted to a
# number) without disturbing the destination cell's existing style. A
# scratch cell is formatted as Text, loaded with the value, then its value
# (not its format) is pasted onto the destination.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

function Set-TextValue($cell, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2 data
Set-TextValue $ws.Range("A2") "1"
$ws.Range("B2").Value = "127.0.0.1"
Set-TextValue $ws.Range("C2") "1488"
$ws.Range("D2").Value = "test"
$ws.Range("E2").Value = "testUser"
$ws.Range("F2").Value = "testPass"

# Row 3 data (same connection details as row 2, different ConnectionID)
Set-TextValue $ws.Range("A3") "2"
$ws.Range("B3").Value = "127.0.0.1"
Set-TextValue $ws.Range("C3") "1488"
$ws.Range("D3").Value = "test"
$ws.Range("E3").Value = "testUser"
$ws.Range("F3").Value = "testPass"

# Clean up scratch cell
$scratch.NumberFormat = "General"
$scratch.Clear()

# --- Formatting updates ---

# Header row (s=3): switch the bold header font to DengXian.
$ws.Range("A1:F1").Font.Name = "DengXian"
$ws.Range("A1:F1").Font.Family = 30

# Row 2 (s=1): white text on a dark grey fill (was blue).
$ws.Range("A2:F2").Font.Family = 24
$ws.Range("A2:F2").Font.Color = 16777215
$ws.Range("A2:F2").Interior.Color = 5855577

# Row 3 (s=2): keep black text/orange fill, just swap the font family.
$ws.Range("A3:F3").Font.Family = 9
